# Refresh the generated financial_summary sheet with the new OCR pass:
#   - the "Current Year"/"Previous Year" labels roll from 2025/2024 to 2024/2023
#   - the Balance Sheet section gains "Property, Plant and Equipment" and
#     "Retained Earnings" lines and the three old "Total Liabilities / Total
#     Equity / Total Liabilities and Equity" rows collapse into a single
#     "Total Equity & Liabilities" row
#   - the Income Statement section gains "General and Administrative
#     Expenses" and "Profit Before Tax" lines
#   - the old "Cash Flow" mini-section (Net Change in Cash / Cash Flow from
#     Operating Activities, row 15) is dropped entirely, and every figure is
#     replaced with the real extracted numbers (vs. the old placeholder
#     500/6,000/.../0 values)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Current/Previous-year figures are comma-formatted numeric-looking text
# ("1,000", "-12,443,892.15", ...) that must stay literal strings, not get
# silently coerced into real numbers by the COM value setter. Pre-marking the
# two value columns as Text handles that; the temporary format is stripped
# again below once every value is in place, so the cells end up with no
# explicit style (matching how this sheet looked before the edit).
$ws.Range("B2:C8").NumberFormat = "@"
$ws.Range("B10:C12").NumberFormat = "@"

# --- Balance Sheet -------------------------------------------------------
$ws.Range("A1").Value = ' Indicator'
$ws.Range("B1").Value = 'Current Year (2024)'
$ws.Range("C1").Value = 'Previous Year (2023) '

$ws.Range("A2").Value = ' Cash and Cash Equivalents'
$ws.Range("B2").Value = '1,000'
$ws.Range("C2").Value = '1,000                '

$ws.Range("A3").Value = ' Accounts Receivable'
$ws.Range("B3").Value = '11,987,605.97'
$ws.Range("C3").Value = '10,711,454.12        '

$ws.Range("A4").Value = ' Property, Plant and Equipment'
$ws.Range("B4").Value = '3,489,523.92'
$ws.Range("C4").Value = '3,494,523.92         '

$ws.Range("A5").Value = ' Total Assets'
$ws.Range("B5").Value = '14,355,193.96'
$ws.Range("C5").Value = '13,424,369.47        '

$ws.Range("A6").Value = ' Accounts Payable'
$ws.Range("B6").Value = '-12,443,892.15'
$ws.Range("C6").Value = '-10,979,515.78       '

$ws.Range("A7").Value = ' Retained Earnings'
$ws.Range("B7").Value = '-2,444,853.69'
$ws.Range("C7").Value = '-2,741,596.38        '

$ws.Range("A8").Value = ' Total Equity & Liabilities'
$ws.Range("B8").Value = '-14,888,745.84'
$ws.Range("C8").Value = '13,721,112.16        '

# --- Income Statement ------------------------------------------------------
$ws.Range("A9").Value = ' Indicator'
$ws.Range("B9").Value = 'Current Year (2024)'
$ws.Range("C9").Value = 'Previous Year (2023) '

$ws.Range("A10").Value = ' Revenue'
$ws.Range("B10").Value = '-1,276,151.85'
$ws.Range("C10").Value = '-1,727,145.61        '

$ws.Range("A11").Value = ' Cost of Goods Sold (COGS)'
$ws.Range("B11").Value = '-367,148.33'
$ws.Range("C11").Value = '-428,513.69          '

$ws.Range("A12").Value = ' General and Administrative Expenses'
$ws.Range("B12").Value = '937,434.64'
$ws.Range("C12").Value = '1,105,786.47         '

$ws.Range("A13").Value = ' Profit Before Tax'
$ws.Range("B13").Value = 'N/A'
$ws.Range("C13").Value = 'N/A                  '

$ws.Range("A14").Value = ' Net Profit'
$ws.Range("B14").Value = 'N/A'
$ws.Range("C14").Value = 'N/A                  '

# Strip the temporary Text format again so B2:C8/B10:C12 are left with no
# explicit style, same as a freshly authored sheet.
$ws.Range("B2:C8").ClearFormats()
$ws.Range("B10:C12").ClearFormats()

# The old "Cash Flow" mini-section (row 15) no longer exists in the rebuilt
# table, so the used range shrinks from A1:C15 to A1:C14.
$ws.Rows.Item(15).Delete()

Write-Output "financial_summary sheet refreshed"
